$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Task 11 ("Crear presentación") is renamed/expanded to mention the defense too.
$ws.Range("B12").Value = "Crear presentación y preparar defensa"

# The "(ongoing)" placeholder comments in the TIEMPO REAL column are replaced
# with actual logged hours now that the tasks are finished.
$ws.Range("E8").Value = 10
$ws.Range("E10").Value = 3
$ws.Range("E11").Value = 5
$ws.Range("E11").HorizontalAlignment = -4108
$ws.Range("E12").Value = 6

# Move the active selection to reflect where the author was last working.
$ws.Range("E11").Select()
